# average with safety stocks
#
# Productdata sheet: InventoryCosts (D), BackorderCosts (F) and LostSale (I)
# are rescaled down by a factor of 0.0004 (i.e. /2500), keeping the
# internal ratios F = 2*D and I = 20*D that already held before the edit.
#
# ForcastedStandardDeviation sheet: the safety-stock standard deviations
# for products 1-4 (columns B-E) in time buckets 7-9 (rows 9-11) are
# zeroed out, matching the other already-zero buckets/products.

$wb = $excel.ActiveWorkbook

$productData = $wb.Worksheets.Item("Productdata")
for ($r = 2; $r -le 11; $r++) {
    $oldD = $productData.Cells.Item($r, 4).Value()
    $newD = $oldD * 0.0004
    $productData.Cells.Item($r, 4).Value = $newD
    $productData.Cells.Item($r, 6).Value = $newD * 2
    $productData.Cells.Item($r, 9).Value = $newD * 20
}

$stdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")
for ($r = 9; $r -le 11; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $stdDev.Cells.Item($r, $c).Value = 0
    }
}
